$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("B1").Value = "FALSE_count"
$ws.Range("C1").Value = "FALSE_percent"
$ws.Range("D1").Value = "TRUE_count"
$ws.Range("E1").Value = "TRUE_percent"

# Update row 2 values
$ws.Range("C2").Value = 5.555555555555555
$ws.Range("E2").Value = 2.702702702702703

# Update row 3 values (B and D swapped, C and E recalculated)
$ws.Range("B3").Value = 17
$ws.Range("C3").Value = 94.44444444444444
$ws.Range("D3").Value = 72
$ws.Range("E3").Value = 97.29729729729731
